$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the vendor ("Leverancier") name columns more explicit
$ws.Range("B1").Value = "LeverancierNaam0001"
$ws.Range("B2").Value = "LeverancierNaam0002"

# Scroll the view back to the top-left corner (topLeftCell A1) and move
# the active selection to E1, matching the saved view state.
$excel.Goto($ws.Range("A1"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E1").Select()
